# Update NATMI ligand-receptor TPM-derived statistics for Col2a1-Ddr1.
# The underlying TPM recount changed the "Ligand-expressing cells" count
# for the ECs cluster (1 -> 2) and the per-cluster total expression values
# (ligand total expression for ECs/MuSCs, receptor total expression for
# ECs/MuSCs target clusters). This cascades into every derived
# detection-rate / average-expression / specificity / edge-weight column.
# The values below are the recomputed outputs written back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Col2a1 -> Ddr1 -> ECs
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.009821333333333333
$ws.Range("H2").Value = 0.029464
$ws.Range("I2").Value = 0.06297798848338983
$ws.Range("J2").Value = 0.06297798848338984
$ws.Range("M2").Value = 0.141694
$ws.Range("N2").Value = 0.425082
$ws.Range("O2").Value = 0.01763793963212447
$ws.Range("P2").Value = 0.01763793963212447
$ws.Range("Q2").Value = 0.001391624005333333
$ws.Range("R2").Value = 0.012524616048
$ws.Range("S2").Value = 0.00111080195902266
$ws.Range("T2").Value = 0.00111080195902266

# Row 3: ECs -> Col2a1 -> Ddr1 -> FAPs
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.009821333333333333
$ws.Range("H3").Value = 0.029464
$ws.Range("I3").Value = 0.06297798848338983
$ws.Range("J3").Value = 0.06297798848338984
$ws.Range("O3").Value = 0.2714637835982539
$ws.Range("P3").Value = 0.2714637835982538
$ws.Range("Q3").Value = 0.02141834736444444
$ws.Range("R3").Value = 0.19276512628
$ws.Range("S3").Value = 0.01709624303710826
$ws.Range("T3").Value = 0.01709624303710826

# Row 4: ECs -> Col2a1 -> Ddr1 -> MuSCs
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.009821333333333333
$ws.Range("H4").Value = 0.029464
$ws.Range("I4").Value = 0.06297798848338983
$ws.Range("J4").Value = 0.06297798848338984
$ws.Range("M4").Value = 5.710985666666667
$ws.Range("N4").Value = 17.132957
$ws.Range("O4").Value = 0.7108982767696218
$ws.Range("P4").Value = 0.7108982767696217
$ws.Range("Q4").Value = 0.05608949389422222
$ws.Range("R4").Value = 0.5048054450480001
$ws.Range("S4").Value = 0.04477094348725891
$ws.Range("T4").Value = 0.04477094348725892

# Row 5: FAPs -> Col2a1 -> Ddr1 -> ECs
$ws.Range("I5").Value = 0.3247949111459754
$ws.Range("J5").Value = 0.3247949111459754
$ws.Range("M5").Value = 0.141694
$ws.Range("N5").Value = 0.425082
$ws.Range("O5").Value = 0.01763793963212447
$ws.Range("P5").Value = 0.01763793963212447
$ws.Range("Q5").Value = 0.007176990025333333
$ws.Range("R5").Value = 0.064592910228
$ws.Range("S5").Value = 0.005728713035613944
$ws.Range("T5").Value = 0.005728713035613944

# Row 6: FAPs -> Col2a1 -> Ddr1 -> FAPs
$ws.Range("I6").Value = 0.3247949111459754
$ws.Range("J6").Value = 0.3247949111459754
$ws.Range("O6").Value = 0.2714637835982539
$ws.Range("P6").Value = 0.2714637835982538
$ws.Range("S6").Value = 0.08817005547314516
$ws.Range("T6").Value = 0.08817005547314516

# Row 7: FAPs -> Col2a1 -> Ddr1 -> MuSCs
$ws.Range("I7").Value = 0.3247949111459754
$ws.Range("J7").Value = 0.3247949111459754
$ws.Range("M7").Value = 5.710985666666667
$ws.Range("N7").Value = 17.132957
$ws.Range("O7").Value = 0.7108982767696218
$ws.Range("P7").Value = 0.7108982767696217
$ws.Range("Q7").Value = 0.2892690386642222
$ws.Range("R7").Value = 2.603421347978
$ws.Range("S7").Value = 0.2308961426372163
$ws.Range("T7").Value = 0.2308961426372163

# Row 8: MuSCs -> Col2a1 -> Ddr1 -> ECs
$ws.Range("G8").Value = 0.09547600000000001
$ws.Range("H8").Value = 0.286428
$ws.Range("I8").Value = 0.6122271003706348
$ws.Range("J8").Value = 0.6122271003706349
$ws.Range("M8").Value = 0.141694
$ws.Range("N8").Value = 0.425082
$ws.Range("O8").Value = 0.01763793963212447
$ws.Range("P8").Value = 0.01763793963212447
$ws.Range("Q8").Value = 0.013528376344
$ws.Range("R8").Value = 0.121755387096
$ws.Range("S8").Value = 0.01079842463748786
$ws.Range("T8").Value = 0.01079842463748786

# Row 9: MuSCs -> Col2a1 -> Ddr1 -> FAPs
$ws.Range("G9").Value = 0.09547600000000001
$ws.Range("H9").Value = 0.286428
$ws.Range("I9").Value = 0.6122271003706348
$ws.Range("J9").Value = 0.6122271003706349
$ws.Range("O9").Value = 0.2714637835982539
$ws.Range("P9").Value = 0.2714637835982538
$ws.Range("Q9").Value = 0.2082139016733333
$ws.Range("R9").Value = 1.87392511506
$ws.Range("S9").Value = 0.1661974850880005
$ws.Range("T9").Value = 0.1661974850880005

# Row 10: MuSCs -> Col2a1 -> Ddr1 -> MuSCs
$ws.Range("G10").Value = 0.09547600000000001
$ws.Range("H10").Value = 0.286428
$ws.Range("I10").Value = 0.6122271003706348
$ws.Range("J10").Value = 0.6122271003706349
$ws.Range("M10").Value = 5.710985666666667
$ws.Range("N10").Value = 17.132957
$ws.Range("O10").Value = 0.7108982767696218
$ws.Range("P10").Value = 0.7108982767696217
$ws.Range("Q10").Value = 0.5452620675106667
$ws.Range("R10").Value = 4.907358607596001
$ws.Range("S10").Value = 0.4352311906451465
$ws.Range("T10").Value = 0.4352311906451465
